$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.2160193333333333
$ws.Range("H2").Value = 0.648058
$ws.Range("I2").Value = 0.02486881244588016
$ws.Range("J2").Value = 0.02486881244588016
$ws.Range("M2").Value = 1.646588666666666
$ws.Range("N2").Value = 4.939766
$ws.Range("O2").Value = 0.039310317935267
$ws.Range("P2").Value = 0.039310317935267
$ws.Range("Q2").Value = 0.3556949860475555
$ws.Range("R2").Value = 3.201254874428
$ws.Range("S2").Value = 0.0009776009239200738
$ws.Range("T2").Value = 0.0009776009239200743

$ws.Range("G3").Value = 0.2160193333333333
$ws.Range("H3").Value = 0.648058
$ws.Range("I3").Value = 0.02486881244588016
$ws.Range("J3").Value = 0.02486881244588016
$ws.Range("O3").Value = 0.278787195370394
$ws.Range("P3").Value = 0.278787195370394
$ws.Range("Q3").Value = 2.522574549786222
$ws.Range("R3").Value = 22.703170948076
$ws.Range("S3").Value = 0.006933106473979276
$ws.Range("T3").Value = 0.006933106473979278

$ws.Range("G4").Value = 0.2160193333333333
$ws.Range("H4").Value = 0.648058
$ws.Range("I4").Value = 0.02486881244588016
$ws.Range("J4").Value = 0.02486881244588016
$ws.Range("M4").Value = 0.7553226666666667
$ws.Range("N4").Value = 2.265968
$ws.Range("O4").Value = 0.01803241742850595
$ws.Range("P4").Value = 0.01803241742850595
$ws.Range("Q4").Value = 0.1631642989048889
$ws.Range("R4").Value = 1.468478690144
$ws.Range("S4").Value = 0.000448444806975335
$ws.Range("T4").Value = 0.0004484448069753351

$ws.Range("G5").Value = 0.2160193333333333
$ws.Range("H5").Value = 0.648058
$ws.Range("I5").Value = 0.02486881244588016
$ws.Range("J5").Value = 0.02486881244588016
$ws.Range("M5").Value = 27.21325766666666
$ws.Range("N5").Value = 81.63977299999999
$ws.Range("O5").Value = 0.6496836961088899
$ws.Range("P5").Value = 0.6496836961088899
$ws.Range("Q5").Value = 5.878589778981556
$ws.Range("R5").Value = 52.90730801083399
$ws.Range("S5").Value = 0.01615686198767818
$ws.Range("T5").Value = 0.01615686198767819

$ws.Range("G6").Value = 0.2160193333333333
$ws.Range("H6").Value = 0.648058
$ws.Range("I6").Value = 0.02486881244588016
$ws.Range("J6").Value = 0.02486881244588016
$ws.Range("M6").Value = 0.5942236666666667
$ws.Range("N6").Value = 1.782671
$ws.Range("O6").Value = 0.01418637315694314
$ws.Range("P6").Value = 0.01418637315694314
$ws.Range("Q6").Value = 0.1283638003242222
$ws.Range("R6").Value = 1.155274202918
$ws.Range("S6").Value = 0.0003527982533272877
$ws.Range("T6").Value = 0.0003527982533272878

$ws.Range("I7").Value = 0.9551554900377276
$ws.Range("J7").Value = 0.9551554900377278
$ws.Range("M7").Value = 1.646588666666666
$ws.Range("N7").Value = 4.939766
$ws.Range("O7").Value = 0.039310317935267
$ws.Range("P7").Value = 0.039310317935267
$ws.Range("Q7").Value = 13.66144923251044
$ws.Range("R7").Value = 122.953043092594
$ws.Range("S7").Value = 0.03754746599099883
$ws.Range("T7").Value = 0.03754746599099883

$ws.Range("I8").Value = 0.9551554900377276
$ws.Range("J8").Value = 0.9551554900377278
$ws.Range("O8").Value = 0.278787195370394
$ws.Range("P8").Value = 0.278787195370394
$ws.Range("S8").Value = 0.2662851202102524
$ws.Range("T8").Value = 0.2662851202102524

$ws.Range("I9").Value = 0.9551554900377276
$ws.Range("J9").Value = 0.9551554900377278
$ws.Range("M9").Value = 0.7553226666666667
$ws.Range("N9").Value = 2.265968
$ws.Range("O9").Value = 0.01803241742850595
$ws.Range("P9").Value = 0.01803241742850595
$ws.Range("Q9").Value = 6.266775955479111
$ws.Range("R9").Value = 56.40098359931199
$ws.Range("S9").Value = 0.01722376250548946
$ws.Range("T9").Value = 0.01722376250548946

$ws.Range("I10").Value = 0.9551554900377276
$ws.Range("J10").Value = 0.9551554900377278
$ws.Range("M10").Value = 27.21325766666666
$ws.Range("N10").Value = 81.63977299999999
$ws.Range("O10").Value = 0.6496836961088899
$ws.Range("P10").Value = 0.6496836961088899
$ws.Range("Q10").Value = 225.7834914028674
$ws.Range("R10").Value = 2032.051422625806
$ws.Range("S10").Value = 0.6205489491264088
$ws.Range("T10").Value = 0.6205489491264089

$ws.Range("I11").Value = 0.9551554900377276
$ws.Range("J11").Value = 0.9551554900377278
$ws.Range("M11").Value = 0.5942236666666667
$ws.Range("N11").Value = 1.782671
$ws.Range("O11").Value = 0.01418637315694314
$ws.Range("P11").Value = 0.01418637315694314
$ws.Range("Q11").Value = 4.930166603998778
$ws.Range("R11").Value = 44.37149943598899
$ws.Range("S11").Value = 0.01355019220457809
$ws.Range("T11").Value = 0.01355019220457809

$ws.Range("G12").Value = 0.173491
$ws.Range("H12").Value = 0.520473
$ws.Range("I12").Value = 0.01997281943922393
$ws.Range("J12").Value = 0.01997281943922393
$ws.Range("M12").Value = 1.646588666666666
$ws.Range("N12").Value = 4.939766
$ws.Range("O12").Value = 0.039310317935267
$ws.Range("P12").Value = 0.039310317935267
$ws.Range("Q12").Value = 0.2856683143686666
$ws.Range("R12").Value = 2.571014829318
$ws.Range("S12").Value = 0.0007851378822195737
$ws.Range("T12").Value = 0.000785137882219574

$ws.Range("G13").Value = 0.173491
$ws.Range("H13").Value = 0.520473
$ws.Range("I13").Value = 0.01997281943922393
$ws.Range("J13").Value = 0.01997281943922393
$ws.Range("O13").Value = 0.278787195370394
$ws.Range("P13").Value = 0.278787195370394
$ws.Range("Q13").Value = 2.025948207800667
$ws.Range("R13").Value = 18.233533870206
$ws.Range("S13").Value = 0.005568166315100524
$ws.Range("T13").Value = 0.005568166315100526

$ws.Range("G14").Value = 0.173491
$ws.Range("H14").Value = 0.520473
$ws.Range("I14").Value = 0.01997281943922393
$ws.Range("J14").Value = 0.01997281943922393
$ws.Range("M14").Value = 0.7553226666666667
$ws.Range("N14").Value = 2.265968
$ws.Range("O14").Value = 0.01803241742850595
$ws.Range("P14").Value = 0.01803241742850595
$ws.Range("Q14").Value = 0.1310416847626666
$ws.Range("R14").Value = 1.179375162864
$ws.Range("S14").Value = 0.0003601582173522639
$ws.Range("T14").Value = 0.0003601582173522641

$ws.Range("G15").Value = 0.173491
$ws.Range("H15").Value = 0.520473
$ws.Range("I15").Value = 0.01997281943922393
$ws.Range("J15").Value = 0.01997281943922393
$ws.Range("M15").Value = 27.21325766666666
$ws.Range("N15").Value = 81.63977299999999
$ws.Range("O15").Value = 0.6496836961088899
$ws.Range("P15").Value = 0.6496836961088899
$ws.Range("Q15").Value = 4.721255285847666
$ws.Range("R15").Value = 42.49129757262899
$ws.Range("S15").Value = 0.01297601515499049
$ws.Range("T15").Value = 0.01297601515499049

$ws.Range("G16").Value = 0.173491
$ws.Range("H16").Value = 0.520473
$ws.Range("I16").Value = 0.01997281943922393
$ws.Range("J16").Value = 0.01997281943922393
$ws.Range("M16").Value = 0.5942236666666667
$ws.Range("N16").Value = 1.782671
$ws.Range("O16").Value = 0.01418637315694314
$ws.Range("P16").Value = 0.01418637315694314
$ws.Range("Q16").Value = 0.1030924581536667
$ws.Range("R16").Value = 0.927832123383
$ws.Range("S16").Value = 0.0002833418695610784
$ws.Range("T16").Value = 0.0002833418695610785

$ws.Range("G17").Value = 0.000025
$ws.Range("H17").Value = 0.00007499999999999999
$ws.Range("I17").Value = 0.000002878077168156263
$ws.Range("J17").Value = 0.000002878077168156264
$ws.Range("M17").Value = 1.646588666666666
$ws.Range("N17").Value = 4.939766
$ws.Range("O17").Value = 0.039310317935267
$ws.Range("P17").Value = 0.039310317935267
$ws.Range("Q17").Value = 0.00004116471666666666
$ws.Range("R17").Value = 0.0003704824499999999
$ws.Range("S17").Value = 0.0000001131381285224556
$ws.Range("T17").Value = 0.0000001131381285224556

$ws.Range("G18").Value = 0.000025
$ws.Range("H18").Value = 0.00007499999999999999
$ws.Range("I18").Value = 0.000002878077168156263
$ws.Range("J18").Value = 0.000002878077168156264
$ws.Range("O18").Value = 0.278787195370394
$ws.Range("P18").Value = 0.278787195370394
$ws.Range("Q18").Value = 0.0002919385166666666
$ws.Range("R18").Value = 0.00262744665
$ws.Range("S18").Value = 0.0000008023710617698503
$ws.Range("T18").Value = 0.0000008023710617698506

$ws.Range("G19").Value = 0.000025
$ws.Range("H19").Value = 0.00007499999999999999
$ws.Range("I19").Value = 0.000002878077168156263
$ws.Range("J19").Value = 0.000002878077168156264
$ws.Range("M19").Value = 0.7553226666666667
$ws.Range("N19").Value = 2.265968
$ws.Range("O19").Value = 0.01803241742850595
$ws.Range("P19").Value = 0.01803241742850595
$ws.Range("Q19").Value = 0.00001888306666666667
$ws.Range("R19").Value = 0.0001699476
$ws.Range("S19").Value = 0.00000005189868888764604
$ws.Range("T19").Value = 0.00000005189868888764605

$ws.Range("G20").Value = 0.000025
$ws.Range("H20").Value = 0.00007499999999999999
$ws.Range("I20").Value = 0.000002878077168156263
$ws.Range("J20").Value = 0.000002878077168156264
$ws.Range("M20").Value = 27.21325766666666
$ws.Range("N20").Value = 81.63977299999999
$ws.Range("O20").Value = 0.6496836961088899
$ws.Range("P20").Value = 0.6496836961088899
$ws.Range("Q20").Value = 0.0006803314416666665
$ws.Range("R20").Value = 0.006122982974999999
$ws.Range("S20").Value = 0.000001869839812294368
$ws.Range("T20").Value = 0.000001869839812294368

$ws.Range("G21").Value = 0.000025
$ws.Range("H21").Value = 0.00007499999999999999
$ws.Range("I21").Value = 0.000002878077168156263
$ws.Range("J21").Value = 0.000002878077168156264
$ws.Range("M21").Value = 0.5942236666666667
$ws.Range("N21").Value = 1.782671
$ws.Range("O21").Value = 0.01418637315694314
$ws.Range("P21").Value = 0.01418637315694314
$ws.Range("Q21").Value = 0.00001485559166666667
$ws.Range("R21").Value = 0.000133700325
$ws.Range("S21").Value = 0.00000004082947668194293
$ws.Range("T21").Value = 0.00000004082947668194294
